$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": rename/retype existing fields, insert two new field rows
# for "new_table" (effect, unique_id) before the trailing blank row.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2: id -> subject_id
$wsOverview.Range("B2").Value = "subject_id"

# Row 3: drug -> entry_date (varchar -> date, max length 12 -> 10)
$wsOverview.Range("B3").Value = "entry_date"
$wsOverview.Range("C3").Value = "date"
$wsOverview.Range("D3").Value = 10

# Row 4: treatment -> discharge_date (varchar -> date, max length 8 -> 0, fraction empty 0 -> 1)
$wsOverview.Range("B4").Value = "discharge_date"
$wsOverview.Range("C4").Value = "date"
$wsOverview.Range("D4").Value = 0
$wsOverview.Range("G4").Value = 1

# Row 5: date -> drug (type/length/etc. unchanged)
$wsOverview.Range("B5").Value = "drug"

# Insert two fresh rows ahead of the old trailing blank row (currently row 6),
# pushing it down to row 8, and populate the two new rows with the new fields.
$wsOverview.Rows.Item(6).Insert()
$wsOverview.Rows.Item(6).Insert()

$wsOverview.Range("A6").Value = "new_table"
$wsOverview.Range("B6").Value = "effect"
$wsOverview.Range("C6").Value = "varchar"
$wsOverview.Range("D6").Value = 14
$wsOverview.Range("E6").Value = 2
$wsOverview.Range("F6").Value = 2
$wsOverview.Range("G6").Value = 0

$wsOverview.Range("A7").Value = "new_table"
$wsOverview.Range("B7").Value = "unique_id"
$wsOverview.Range("C7").Value = "int"
$wsOverview.Range("D7").Value = 1
$wsOverview.Range("E7").Value = 2
$wsOverview.Range("F7").Value = 2
$wsOverview.Range("G7").Value = 0

# ---------------------------------------------------------------------------
# Sheet "new_table": the frequency-summary sheet. Relabel the existing
# id/drug/treatment/date column pairs to match the renamed fields, and add
# two more column pairs for the new "effect" / "unique_id" fields.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("new_table")

$wsSummary.Range("A1").Value = "subject_id"
$wsSummary.Range("C1").Value = "entry_date"
$wsSummary.Range("E1").Value = "discharge_date"
$wsSummary.Range("G1").Value = "drug"

$wsSummary.Range("I1").Value = "effect"
$wsSummary.Range("J1").Value = "Frequency"
$wsSummary.Range("K1").Value = "unique_id"
$wsSummary.Range("L1").Value = "Frequency"

$wsSummary.Range("I2").Value = "List truncated..."
$wsSummary.Range("K2").Value = "List truncated..."
# J2 / L2 stay blank, mirroring the existing B2/D2/F2/H2 "Frequency" placeholder
# columns, which are likewise empty beneath their header.
